$d = $word.ActiveDocument

# wdLineSpaceMultiple = 5; 276/240 = 1.15 multiple -> LineSpacing expressed as
# 1.15 * 12pt single-line-height = 13.8
$wdLineSpaceMultiple = 5
$lineSpacingValue = 13.8

# --- Paragraph 1 ("CHAPTER 3") and Paragraph 2 ("WIRELESS PACKET ANALYSIS"):
# add 1.15 line spacing to their paragraph formatting.
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LineSpacingRule = $wdLineSpaceMultiple
$p1.Range.ParagraphFormat.LineSpacing = $lineSpacingValue

$p2 = $d.Paragraphs(2)
$p2.Range.ParagraphFormat.LineSpacingRule = $wdLineSpaceMultiple
$p2.Range.ParagraphFormat.LineSpacing = $lineSpacingValue

# --- Insert a new empty paragraph right before the trailing (bookmarked) empty
# paragraph, carrying the same run formatting (Times New Roman, sz 24/szCs 30)
# as that paragraph, then give both that new paragraph and the original
# bookmarked paragraph the same 1.15 line spacing.
$p3 = $d.Paragraphs(3)
$p3.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(3)
$bookmarkPara = $d.Paragraphs(4)

$newPara.Range.ParagraphFormat.LineSpacingRule = $wdLineSpaceMultiple
$newPara.Range.ParagraphFormat.LineSpacing = $lineSpacingValue

$bookmarkPara.Range.ParagraphFormat.LineSpacingRule = $wdLineSpaceMultiple
$bookmarkPara.Range.ParagraphFormat.LineSpacing = $lineSpacingValue
